# Append one new row (row 87) to Sheet1 with the new game entry, matching
# the existing un-styled data rows (A2:D86) -- no fill/border/alignment,
# just plain text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the value "666" which looks numeric; force it to be
# stored as text (like every other cell in the sheet) instead of a
# number, then drop the temporary number-format override so the cell
# keeps the same (default/no) style as its neighbours.
$ws.Range("A87").NumberFormat = "@"
$ws.Range("A87").Value = "666"
$ws.Range("A87").Style = "Normal"

$ws.Range("B87").Value = "Incompleto"
$ws.Range("C87").Value = "PC"
$ws.Range("D87").Value = "Platinado"
